$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 267.863445493107
$ws.Range("G2").Value = 20.02461998266713
$ws.Range("H2").Value = 555.9628879047845
$ws.Range("I2").Value = 1.426082368387081
$ws.Range("J2").Value = 0.5264331343036044
$ws.Range("K2").Value = 2.54411291104781
$ws.Range("L2").Value = 0.2153346801364407
$ws.Range("M2").Value = 0.03228456974986398
$ws.Range("N2").Value = 0.4442538880519111

$ws.Range("F3").Value = 0.02025071042173724
$ws.Range("G3").Value = 0.01354245178581307
$ws.Range("H3").Value = 0.02745494690914944
$ws.Range("I3").Value = 0.01886812045533984
$ws.Range("J3").Value = 0.0125218690743075
$ws.Range("K3").Value = 0.02564939419428678
$ws.Range("L3").Value = 0.02016558230488541
$ws.Range("M3").Value = 0.01346568465930349
$ws.Range("N3").Value = 0.02736159213845523

$ws.Range("F4").Value = 267.8836962035288
$ws.Range("G4").Value = 20.03816243445294
$ws.Range("H4").Value = 555.9903428516936
$ws.Range("I4").Value = 1.444950488842421
$ws.Range("J4").Value = 0.5389550033779119
$ws.Range("K4").Value = 2.569762305242096
$ws.Range("L4").Value = 0.2355002624413261
$ws.Range("M4").Value = 0.04575025440916747
$ws.Range("N4").Value = 0.4716154801903665
